# Apply the "inverse availability, support dates, support results upload" edit.
#
# Summary of the functional change (org_data sheet):
#   - H2/J2/L2 and H3/J3/L3 used to hold bare day-of-month integers
#     (first_date/second_date/third_date). They now hold real Excel date
#     serial values, formatted with a date number format (d-mmm), matching
#     the actual dates those day numbers referred to (October 2023).
#   - The active selection on the org_data sheet moves from F12 to M10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("org_data")

# New cell style: a date number format (built-in numFmtId 16 == "d-mmm"),
# applied the same way Excel applies a number format to a range.
$ws.Range("H2").NumberFormat = "d-mmm"
$ws.Range("J2").NumberFormat = "d-mmm"
$ws.Range("L2").NumberFormat = "d-mmm"
$ws.Range("H3").NumberFormat = "d-mmm"
$ws.Range("J3").NumberFormat = "d-mmm"
$ws.Range("L3").NumberFormat = "d-mmm"

# Row 2 dates: first_date, second_date, third_date
$ws.Range("H2").Value = 45205
$ws.Range("J2").Value = 45218
$ws.Range("L2").Value = 45228

# Row 3 dates: first_date, second_date, third_date
$ws.Range("H3").Value = 45209
$ws.Range("J3").Value = 45215
$ws.Range("L3").Value = 45217

# Move the sheet's active selection from F12 to M10.
$ws.Range("M10").Select()

$wb.Save()
